$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 30000
$ws.Range("J3").Value = 30000
$ws.Range("L3").Value = 30000
$ws.Range("N3").Value = -30228
$ws.Range("H17").Value = 84347.36
$ws.Range("J17").Value = 84347.36
$ws.Range("L17").Value = 253042.08
$ws.Range("N17").Value = -253378.08
$ws.Range("H97").Value = 25610.125
$ws.Range("J97").Value = 25610.125
$ws.Range("L97").Value = 76830.375
$ws.Range("N97").Value = -77822.375
$ws.Range("H99").Value = 1436.6428
$ws.Range("I99").Value = 1194
$ws.Range("J99").Value = 2326.3333
$ws.Range("K99").Value = 3582
$ws.Range("L99").Value = 6978.999899999999
$ws.Range("M99").Value = -2084
$ws.Range("N99").Value = -9974.999899999999
$ws.Range("H100").Value = 2106.125
$ws.Range("I100").Value = 975
$ws.Range("J100").Value = 2483.1667
$ws.Range("K100").Value = 975
$ws.Range("L100").Value = 2483.1667
$ws.Range("M100").Value = -434
$ws.Range("N100").Value = -3565.1667
$ws.Range("H101").Value = 566.2
$ws.Range("I101").Value = 413.2857
$ws.Range("J101").Value = 700
$ws.Range("K101").Value = 1239.8571
$ws.Range("L101").Value = 2100
$ws.Range("M101").Value = 382.1428999999998
$ws.Range("N101").Value = -5344
$ws.Range("H102").Value = 30000
$ws.Range("J102").Value = 30000
$ws.Range("L102").Value = 30000
$ws.Range("N102").Value = -36490
$ws.Range("H129").Value = 836.06665
$ws.Range("I129").Value = 442.2
$ws.Range("J129").Value = 1033
$ws.Range("K129").Value = 1326.6
$ws.Range("L129").Value = 3099
$ws.Range("M129").Value = 3673.4
$ws.Range("N129").Value = -13099
$ws.Range("H132").Value = 4074.5107
$ws.Range("I132").Value = 1963.317
$ws.Range("J132").Value = 18501
$ws.Range("K132").Value = 5889.951
$ws.Range("L132").Value = 55503
$ws.Range("M132").Value = -3359.951
$ws.Range("N132").Value = -60563
$ws.Range("H135").Value = 15385132
$ws.Range("I135").Value = 359.75
$ws.Range("J135").Value = 58824490
$ws.Range("K135").Value = 3237.75
$ws.Range("L135").Value = 529420410
$ws.Range("M135").Value = -702.75
$ws.Range("N135").Value = -529425480
$ws.Range("H137").Value = 958902.5600000001
$ws.Range("I137").Value = 1105.8975
$ws.Range("J137").Value = 2924906.2
$ws.Range("K137").Value = 3317.6925
$ws.Range("L137").Value = 8774718.600000001
$ws.Range("M137").Value = -767.6925000000001
$ws.Range("N137").Value = -8779818.600000001
$ws.Range("H138").Value = 3511174
$ws.Range("I138").Value = 1083.6
$ws.Range("J138").Value = 9095409
$ws.Range("K138").Value = 3250.8
$ws.Range("L138").Value = 27286227
$ws.Range("M138").Value = 1889.2
$ws.Range("N138").Value = -27296507
$ws.Range("H141").Value = 939.7406999999999
$ws.Range("I141").Value = 723.8200000000001
$ws.Range("J141").Value = 3638.75
$ws.Range("K141").Value = 2171.46
$ws.Range("L141").Value = 10916.25
$ws.Range("M141").Value = 3008.54
$ws.Range("N141").Value = -21276.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1243.87
$ws.Range("I32").Value = 1205.7878
$ws.Range("J32").Value = 5014
$ws.Range("K32").Value = 1205.7878
$ws.Range("L32").Value = 5014
$ws.Range("M32").Value = -918.7878000000001
$ws.Range("N32").Value = -5588
$ws.Range("H61").Value = 1050.6405
$ws.Range("I61").Value = 904.3026
$ws.Range("J61").Value = 1906.1538
$ws.Range("K61").Value = 904.3026
$ws.Range("L61").Value = 1906.1538
$ws.Range("M61").Value = -692.3026
$ws.Range("N61").Value = -2330.1538
$ws.Range("H74").Value = 20635.842
$ws.Range("I74").Value = 26654.486
$ws.Range("J74").Value = 7595.4443
$ws.Range("K74").Value = 26654.486
$ws.Range("L74").Value = 7595.4443
$ws.Range("M74").Value = -25780.486
$ws.Range("N74").Value = -9343.444299999999
$ws.Range("H77").Value = 20635.842
$ws.Range("I77").Value = 26654.486
$ws.Range("J77").Value = 7595.4443
$ws.Range("K77").Value = 133272.43
$ws.Range("L77").Value = 37977.2215
$ws.Range("M77").Value = -128904.43
$ws.Range("N77").Value = -46713.2215
$ws.Range("H136").Value = 1050.6405
$ws.Range("I136").Value = 904.3026
$ws.Range("J136").Value = 1906.1538
$ws.Range("K136").Value = 2712.9078
$ws.Range("L136").Value = 5718.4614
$ws.Range("M136").Value = -162.9078
$ws.Range("N136").Value = -10818.4614

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H81").Value = 31820
$ws.Range("J81").Value = 31820
$ws.Range("L81").Value = 31820
$ws.Range("N81").Value = -33942
$ws.Range("H84").Value = 31820
$ws.Range("J84").Value = 31820
$ws.Range("L84").Value = 95460
$ws.Range("N84").Value = -106068
$ws.Range("H134").Value = 479202.5
$ws.Range("I134").Value = 743701.4399999999
$ws.Range("J134").Value = 3104.4
$ws.Range("K134").Value = 2231104.32
$ws.Range("L134").Value = 9313.200000000001
$ws.Range("M134").Value = -2228569.32
$ws.Range("N134").Value = -14383.2
$ws.Range("H140").Value = 20000
$ws.Range("J140").Value = 20000
$ws.Range("L140").Value = 20000
$ws.Range("N140").Value = -30360
$ws.Range("H141").Value = 91000
$ws.Range("J141").Value = 91000
$ws.Range("L141").Value = 91000
$ws.Range("N141").Value = -101360

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 14286902
$ws.Range("I31").Value = 1189.8788
$ws.Range("K31").Value = 1189.8788
$ws.Range("M31").Value = -894.8788
$ws.Range("H34").Value = 14286902
$ws.Range("I34").Value = 1189.8788
$ws.Range("K34").Value = 1189.8788
$ws.Range("M34").Value = -987.8788
$ws.Range("H58").Value = 3524.4358
$ws.Range("I58").Value = 3712.9429
$ws.Range("J58").Value = 1875
$ws.Range("K58").Value = 3712.9429
$ws.Range("L58").Value = 1875
$ws.Range("M58").Value = -3509.9429
$ws.Range("N58").Value = -2281
$ws.Range("H132").Value = 571236.0600000001
$ws.Range("I132").Value = 1465.1052
$ws.Range("J132").Value = 4630854.5
$ws.Range("K132").Value = 4395.3156
$ws.Range("L132").Value = 13892563.5
$ws.Range("M132").Value = -1865.3156
$ws.Range("N132").Value = -13897623.5
$ws.Range("H134").Value = 1714.8948
$ws.Range("I134").Value = 1812.017
$ws.Range("K134").Value = 5436.051
$ws.Range("M134").Value = -2901.051
$ws.Range("H136").Value = 3524.4358
$ws.Range("I136").Value = 3712.9429
$ws.Range("J136").Value = 1875
$ws.Range("K136").Value = 11138.8287
$ws.Range("L136").Value = 5625
$ws.Range("M136").Value = -8588.8287
$ws.Range("N136").Value = -10725

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 1131.9286
$ws.Range("I4").Value = 336.75
$ws.Range("J4").Value = 1450
$ws.Range("K4").Value = 1010.25
$ws.Range("L4").Value = 4350
$ws.Range("M4").Value = -898.25
$ws.Range("N4").Value = -4574
$ws.Range("H80").Value = 2630.6365
$ws.Range("I80").Value = 819.25
$ws.Range("J80").Value = 3665.7144
$ws.Range("K80").Value = 2457.75
$ws.Range("L80").Value = 10997.1432
$ws.Range("M80").Value = -1521.75
$ws.Range("N80").Value = -12869.1432
$ws.Range("H83").Value = 2630.6365
$ws.Range("I83").Value = 819.25
$ws.Range("J83").Value = 3665.7144
$ws.Range("K83").Value = 7373.25
$ws.Range("L83").Value = 32991.4296
$ws.Range("M83").Value = -2693.25
$ws.Range("N83").Value = -42351.4296
$ws.Range("H88").Value = 5363.3335
$ws.Range("J88").Value = 5363.3335
$ws.Range("L88").Value = 16090.0005
$ws.Range("N88").Value = -16946.0005
$ws.Range("H91").Value = 5363.3335
$ws.Range("J91").Value = 5363.3335
$ws.Range("L91").Value = 16090.0005
$ws.Range("N91").Value = -19054.0005
$ws.Range("H109").Value = 2222.762
$ws.Range("I109").Value = 1567.7693
$ws.Range("J109").Value = 3287.125
$ws.Range("K109").Value = 4703.3079
$ws.Range("L109").Value = 9861.375
$ws.Range("M109").Value = -3663.3079
$ws.Range("N109").Value = -11941.375
$ws.Range("H131").Value = 1170.6364
$ws.Range("I131").Value = 950.5
$ws.Range("J131").Value = 1175.1753
$ws.Range("K131").Value = 2851.5
$ws.Range("L131").Value = 3525.525900000001
$ws.Range("M131").Value = 2188.5
$ws.Range("N131").Value = -13605.5259
$ws.Range("H133").Value = 6895.5713
$ws.Range("J133").Value = 7058.8237
$ws.Range("L133").Value = 21176.4711
$ws.Range("N133").Value = -31296.4711

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 1756371.9
$ws.Range("I132").Value = 1920.1052
$ws.Range("J132").Value = 5265275.5
$ws.Range("K132").Value = 5760.3156
$ws.Range("L132").Value = 15795826.5
$ws.Range("M132").Value = -3230.3156
$ws.Range("N132").Value = -15800886.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 3043.44
$ws.Range("I132").Value = 2829.5615
$ws.Range("J132").Value = 3720.7222
$ws.Range("K132").Value = 8488.684499999999
$ws.Range("L132").Value = 11162.1666
$ws.Range("M132").Value = -5958.684499999999
$ws.Range("N132").Value = -16222.1666
$ws.Range("H136").Value = 1538.1786
$ws.Range("I136").Value = 946.3171
$ws.Range("J136").Value = 3155.9333
$ws.Range("K136").Value = 2838.9513
$ws.Range("L136").Value = 9467.7999
$ws.Range("M136").Value = -288.9512999999997
$ws.Range("N136").Value = -14567.7999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 65693
$ws.Range("I122").Value = 10196.1
$ws.Range("J122").Value = 121189.9
$ws.Range("K122").Value = 30588.3
$ws.Range("L122").Value = 363569.7
$ws.Range("M122").Value = -28138.3
$ws.Range("N122").Value = -368469.7
$ws.Range("H132").Value = 2380.843
$ws.Range("I132").Value = 2492.5134
$ws.Range("J132").Value = 2085.7144
$ws.Range("K132").Value = 7477.540199999999
$ws.Range("L132").Value = 6257.1432
$ws.Range("M132").Value = -4947.540199999999
$ws.Range("N132").Value = -11317.1432
$ws.Range("H136").Value = 2057.228
$ws.Range("I136").Value = 1661.5209
$ws.Range("J136").Value = 4167.6665
$ws.Range("K136").Value = 4984.5627
$ws.Range("L136").Value = 12502.9995
$ws.Range("M136").Value = -2434.5627
$ws.Range("N136").Value = -17602.9995
